$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19-107 down to 20-108
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new data record
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C19").Value = "Los Lagos"
$ws.Range("D19").Value = 44859
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 100112031
$ws.Range("G19").Value = "Poroto verde"
$ws.Range("H19").Value = "Magnum"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 35
$ws.Range("K19").Value = 36000
$ws.Range("L19").Value = 36000
$ws.Range("M19").Value = 36000
$ws.Range("N19").Value = "`$/malla 25 kilos"
$ws.Range("O19").Value = "Perú"
$ws.Range("P19").Value = 1440
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
